$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting existing rows 119-144 down to 120-145
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record
$ws.Range("A119").Value = 10
$ws.Range("B119").Value = "Vega Modelo de Temuco"
$ws.Range("C119").Value = "La Araucanía"
$ws.Range("D119").Value = 44522
$ws.Range("E119").Value = 9
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100103
$ws.Range("H119").Value = "Frutos de hueso (carozo)"
$ws.Range("I119").Value = 100103004
$ws.Range("J119").Value = "Durazno"
$ws.Range("K119").Value = "Early Majestic"
$ws.Range("L119").Value = "Primera"
$ws.Range("M119").Value = 80
$ws.Range("N119").Value = 25000
$ws.Range("O119").Value = 25000
$ws.Range("P119").Value = 25000
$ws.Range("Q119").Value = "$/bandeja 15 kilos empedrada"
$ws.Range("R119").Value = "Provincia de Limarí"
$ws.Range("S119").Value = 1667
$ws.Range("T119").Value = 15
